# "able to add trait folders now"
# Adds two new columns to the "Adobe AAM" sheet:
#   - "Trait Folder Path" (inserted right after "Segment Status")
#   - "Data Feed Description" (inserted right after "Data Source Name")
# and populates two sample trait-folder rows (rows 3 and 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Adobe AAM")

$xlPasteFormats = -4122

# --- Insert "Trait Folder Path" column before old column E (Data Source ID) ---
$ws.Columns("E:E").Insert()
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E1").Value = "Trait Folder Path"

$ws.Range("B2").Copy() | Out-Null
$ws.Range("E2").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E2").Value = "Add: Required"

# --- Insert "Data Feed Description" column after "Data Source Name" (now column G) ---
$ws.Columns("H:H").Insert()
$ws.Range("D1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("H1").Value = "Data Feed Description"

$ws.Range("A2").Copy() | Out-Null
$ws.Range("H2").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("H2").Value = "Add: Optional if data source exists"

# --- Row 2 text updates. J2 keeps its "optional" cell look (style of K2) even though the
#     wording now reads "Not Required (default FIXED)"; L2/N2 stay in the same optional style. ---
$ws.Range("K2").Copy() | Out-Null
$ws.Range("J2").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("J2").Value = "Not Required (default FIXED)"

$ws.Range("L2").Value = "Add: Optional (FIXED or CPM)"
$ws.Range("N2").Value = "Add: Optional (FIXED or CPM)"

# --- Row 3: sample trait-folder-path row ---
$ws.Range("E3").Value = "/All Traits/TEST20181029"
$ws.Range("H3").Value = "Test on 20181028"
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = "CPM"

# --- Row 4: second sample trait-folder-path row (sub-folder) ---
$ws.Range("E4").Value = "/All Traits/TEST20181029/TEST"
$ws.Range("G4").Value = "test20181028"

# --- Resize the two new columns to fit their new content ---
$ws.Columns("E:E").ColumnWidth = $ws.Columns("D:D").ColumnWidth
$ws.Columns("H:H").AutoFit() | Out-Null
